# Adds the "k from 5 to 250 (increments of 5)" KNN sweep table for
# Raw/PCA/FLD, plus the printed-array strings for the FLD and PCA
# 10-fold-CV accuracy vectors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: column headers for k = 5, 10, 15, 20
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 20

# Row 20: Raw / KNN accuracies
$ws.Range("B20").Value = "Raw"
$ws.Range("C20").Value = "KNN"
$ws.Range("D20").Value = 0.8554
$ws.Range("E20").Value = 0.8515
$ws.Range("F20").Value = 0.8462
$ws.Range("G20").Value = 0.8415

# Row 21: PCA / KNN accuracies
$ws.Range("B21").Value = "PCA"
$ws.Range("C21").Value = "KNN"
$ws.Range("D21").Value = 0.8603
$ws.Range("E21").Value = 0.8619
$ws.Range("F21").Value = 0.8565
$ws.Range("G21").Value = 0.8541

# I21:I25: full PCA accuracy-vector text, wrapped the way numpy prints it
$ws.Range("I21").Value = "[0.8203 0.8263 0.8273 0.8291 0.8293 0.8295 0.8291 0.8292 0.8288 0.8288"
$ws.Range("I22").Value = " 0.8281 0.8278 0.8281 0.8284 0.8287 0.8281 0.8273 0.8265 0.8271 0.8263"
$ws.Range("I23").Value = " 0.8261 0.8253 0.8254 0.8256 0.8252 0.8256 0.8257 0.826  0.8254 0.825"
$ws.Range("I24").Value = " 0.8249 0.825  0.8246 0.8245 0.8239 0.8232 0.8235 0.8237 0.8236 0.8237"
$ws.Range("I25").Value = " 0.8234 0.8233 0.8235 0.8228 0.8231 0.8233 0.8227 0.8228 0.8227 0.8221]"

# Row 29: FLD / KNN accuracies
$ws.Range("B29").Value = "FLD"
$ws.Range("C29").Value = "KNN"
$ws.Range("D29").Value = 0.8203
$ws.Range("E29").Value = 0.8263
$ws.Range("F29").Value = 0.8273
$ws.Range("G29").Value = 0.8291

# I29:I33: full FLD accuracy-vector text (same underlying shared strings as
# I21:I25 above, since both blocks print the identical numbers)
$ws.Range("I29").Value = "[0.8203 0.8263 0.8273 0.8291 0.8293 0.8295 0.8291 0.8292 0.8288 0.8288"
$ws.Range("I30").Value = " 0.8281 0.8278 0.8281 0.8284 0.8287 0.8281 0.8273 0.8265 0.8271 0.8263"
$ws.Range("I31").Value = " 0.8261 0.8253 0.8254 0.8256 0.8252 0.8256 0.8257 0.826  0.8254 0.825"
$ws.Range("I32").Value = " 0.8249 0.825  0.8246 0.8245 0.8239 0.8232 0.8235 0.8237 0.8236 0.8237"
$ws.Range("I33").Value = " 0.8234 0.8233 0.8235 0.8228 0.8231 0.8233 0.8227 0.8228 0.8227 0.8221]"

# Match the saved selection/cursor position from the source workbook
[void]$ws.Range("I18").Select()
